$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'65.006.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").Formula = "'3.391.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Formula = "'559.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("D6").Formula = "'173.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("E7").Value = "  +1.71%  "

$ws.Range("D8").Formula = "'3.381.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  +11.49%  "

$ws.Range("D11").Formula = "'0.631"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.29%  "

$ws.Range("D12").Formula = "'54.45"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = "  +5.35%  "

$ws.Range("E14").Value = "  +3.09%  "

$ws.Range("D15").Formula = "'3.929.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("D18").Formula = "'3.384.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("D19").Formula = "'64.912.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.10%  "

$ws.Range("D20").Formula = "'11.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("E21").Value = "  +2.15%  "

$ws.Range("D22").Formula = "'471.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.32%  "

$ws.Range("D23").Formula = "'4.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.78%  "

$ws.Range("E24").Value = "  +2.73%  "

$ws.Range("D25").Formula = "'87.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.15%  "

$ws.Range("D26").Formula = "'13.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("E27").Value = "  +7.00%  "

$ws.Range("D28").Formula = "'10.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.51%  "

$ws.Range("E29").Value = "  +2.26%  "

$ws.Range("D30").Formula = "'30.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.85%  "

$ws.Range("D31").Formula = "'6.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.44%  "

$ws.Range("D32").Formula = "'11.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "

$ws.Range("D33").Formula = "'572.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("D34").Formula = "'61.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.73%  "

$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").Formula = "'3.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.33%  "

$ws.Range("D38").Formula = "'0.139"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.90%  "

$ws.Range("D39").Formula = "'35.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.11%  "

$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("D42").Formula = "'3.096.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("D43").Formula = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("E45").Value = "  +3.98%  "

$ws.Range("E46").Value = "  +5.53%  "

$ws.Range("D47").Formula = "'2.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.67%  "

$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").Formula = "'2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Formula = "'139.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.24%  "

$ws.Range("D51").Formula = "'8.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.31%  "
